# Clean up maintenance form
#
# The "survey" sheet (sheet1) gets a new "end screen" / "begin screen"
# pair of clause rows inserted right before the "select_multiple" question
# for common spare parts, and that question's type is switched to
# "select_multiple_inline". The "choices" sheet stays as-is (its row
# values simply track the shared-string table automatically). Finally,
# the "survey" tab becomes the active/selected sheet instead of "choices".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "survey"

# Insert two new blank rows above the old row 9, pushing everything
# (including the "end screen" row that used to be row 11) down by two.
$ws.Rows("9:10").Insert()

# New row 9: clause = "end screen" (closes the previous screen).
$ws.Range("A9").Value = "end screen"

# New row 10: clause = "begin screen" (opens the spare-parts screen).
$ws.Range("A10").Value = "begin screen"

# The spare-parts question (now on row 11 after the insert) changes
# from a plain "select_multiple" to a "select_multiple_inline" widget.
$ws.Range("C11").Value = "select_multiple_inline"

# Make "survey" the active tab/sheet and leave the selection on A11,
# matching where editing last happened.
$ws.Select() | Out-Null
$ws.Range("A11").Select() | Out-Null
